$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Edit Repayment Schedule")

# Insert a new row 6 ("waittopageload1" / 2000), pushing the existing
# rows 6-12 down to 7-13.
$ws.Rows.Item(6).Insert()

$ws.Range("A6").Value = "waittopageload1"
$ws.Range("B6").Value = 2000

# Match the formatting used by the other "amount" cell (B3) exactly by
# copying its format onto the new B6 cell.
$ws.Range("B3").Copy()
$ws.Range("B6").PasteSpecial(-4122)

# Select the new row and make this sheet the active tab, like the
# author did when recording this scenario.
$ws.Activate()
$ws.Range("A6:B6").Select()
